$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "João Martins"
$ws.Range("B7").Value = "GEOGRAFIA"
$ws.Range("A8").Value = "Gabriela Silva"
$ws.Range("B8").Value = "FÍSICA"
$ws.Range("A9").Value = "Ricardo Fabião Amaro"
$ws.Range("B9").Value = "BIOLOGIA"
$ws.Range("A10").Value = "Ricardo Oliveira"
$ws.Range("B10").Value = "QUÍMICA"
$ws.Range("A11").Value = "Ricardo Fabião Amaro"
$ws.Range("B11").Value = "BIOLOGIA"
$ws.Range("A12").Value = "Carlos Lima"
$ws.Range("B12").Value = "HISTÓRIA"
$ws.Range("A13").Value = "Gabriela Silva"
$ws.Range("B13").Value = "FÍSICA"
$ws.Range("A14").Value = "Ricardo Oliveira"
$ws.Range("B14").Value = "QUÍMICA"
$ws.Range("A15").Value = "Manoel Amaro"
$ws.Range("B15").Value = "MATEMÁTICA"
$ws.Range("A16").Value = "Carlos Lima"
$ws.Range("B16").Value = "HISTÓRIA"
$ws.Range("A17").Value = "Gabriela Silva"
$ws.Range("B17").Value = "FÍSICA"
$ws.Range("A18").Value = "Luciana Costa"
$ws.Range("B18").Value = "MATEMÁTICA"
$ws.Range("A19").Value = "Ricardo Fabião Amaro"
$ws.Range("B19").Value = "BIOLOGIA"
$ws.Range("A20").Value = "Manoel Amaro"
$ws.Range("B20").Value = "QUÍMICA"
$ws.Range("A21").Value = "Ricardo Fabião Amaro"
$ws.Range("B21").Value = "BIOLOGIA"
$ws.Range("A22").Value = "João Martins"
$ws.Range("B22").Value = "GEOGRAFIA"
$ws.Range("A23").Value = "Gabriela Silva"
$ws.Range("B23").Value = "FÍSICA"
$ws.Range("A24").Value = "Marcos Pereira"
$ws.Range("B24").Value = "QUÍMICA"
$ws.Range("A25").Value = "Luciana Costa"
$ws.Range("B25").Value = "MATEMÁTICA"
$ws.Range("A26").Value = "João Martins"
$ws.Range("B26").Value = "GEOGRAFIA"
$ws.Range("A27").Value = "Gabriela Silva"
$ws.Range("B27").Value = "FÍSICA"
$ws.Range("A28").Value = "Luciana Costa"
$ws.Range("B28").Value = "MATEMÁTICA"
$ws.Range("A29").Value = "Ricardo Fabião Amaro"
$ws.Range("B29").Value = "BIOLOGIA"
$ws.Range("A30").Value = "Ricardo Oliveira"
$ws.Range("B30").Value = "QUÍMICA"
$ws.Range("A31").Value = "Ricardo Fabião Amaro"
$ws.Range("B31").Value = "BIOLOGIA"
$ws.Range("A32").Value = "Carlos Lima"
$ws.Range("B32").Value = "HISTÓRIA"
$ws.Range("A33").Value = "Gabriela Silva"
$ws.Range("B33").Value = "FÍSICA"
$ws.Range("A34").Value = "Luciana Costa"
$ws.Range("B34").Value = "MATEMÁTICA"
$ws.Range("A38").Value = "Gabriela Silva"
$ws.Range("B38").Value = "FÍSICA"
$ws.Range("A39").Value = "Ricardo Fabião Amaro"
$ws.Range("B39").Value = "BIOLOGIA"
$ws.Range("A40").Value = "Marcos Pereira"
$ws.Range("B40").Value = "QUÍMICA"
$ws.Range("A41").Value = "Ricardo Fabião Amaro"
$ws.Range("B41").Value = "BIOLOGIA"
$ws.Range("B46").Value = "MATEMÁTICA"
$ws.Range("A49").Value = "Ricardo Fabião Amaro"
$ws.Range("B49").Value = "FÍSICA"
$ws.Range("D49").Value = 2
$ws.Range("A51").Value = "Ricardo Fabião Amaro"
$ws.Range("B51").Value = "BIOLOGIA"
$ws.Range("C51").Value = "Sexta"
$ws.Range("D51").Value = 3
